# "Rename example to correct layers"
#
# The shared-string "water board" (used in every data row of the
# "Layer" column, column D) is renamed to "water_board" so the example
# import file matches the underscore-separated layer names the
# importer actually expects.
#
# The accompanying OOXML diff also shows the header row's bordered
# cells (A1:E1) picking up a new solid white fill (fillId referencing
# a freshly-added <fill> in styles.xml) - a cosmetic header highlight
# that Excel added when it resaved the workbook. We reproduce that
# visual effect here too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content fix: water board -> water_board (column D, the "Layer" column) ---
$lastRow = $ws.Cells(1, 1).End(-4121).Row  # xlDown
if ($lastRow -lt 2) { $lastRow = 23 }
$layerRange = $ws.Range("D2:D" + $lastRow)
$layerRange.Value = "water_board"

# --- cosmetic: header row gets a white fill behind its bordered cells ---
$headerRange = $ws.Range("A1:E1")
$headerRange.Interior.ColorIndex = 2
